# Remove the "License, Copyright and Data" slide (show position 15).
# This makes "IP and Licences" its own, separate episode (per commit message),
# and PowerPoint's normal slide-delete behaviour renumbers/shifts the
# remaining slide parts accordingly.
$p = $ppt.ActivePresentation
$p.Slides.Item(15).Delete()
